$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.303.83'
$ws.Range('E2').Value = '  +2.39%  '
$ws.Range('D3').Value = '3.181.79'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'595.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.17%  '
$ws.Range('D6').Value = "'153.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.31%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.181.44'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').Value = "'0.550"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.91%  '
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('E11').Value = '  -2.74%  '
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').Value = "'0.0000268"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'39.00"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.30%  '
$ws.Range('D15').Value = '3.702.17'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').Value = '66.267.86'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D18').Value = '3.183.08'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').Value = "'513.47"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').Value = "'15.34"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('E22').Value = '  +3.31%  '
$ws.Range('D23').Value = "'8.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.39%  '
$ws.Range('D24').Value = "'14.89"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('D28').Value = "'2.99"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.48%  '
$ws.Range('E29').Value = '  +7.41%  '
$ws.Range('D30').Value = "'7.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.85%  '
$ws.Range('D31').Value = "'2.93"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.35%  '
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('D33').Value = "'1.23"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.82%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'6.50"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').Value = "'503.39"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('D37').Value = "'54.81"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('D38').Value = "'0.0899"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('E39').Value = '  +2.31%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = "'0.125"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.88%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').Value = "'8.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.56%  '
$ws.Range('D42').Value = "'0.303"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.04%  '
$ws.Range('E43').Value = '  -1.88%  '
$ws.Range('D44').Value = '0.0₃0675'
$ws.Range('E44').Value = '  +16.20%  '
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').Value = '2.901.99'
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').Value = "'28.55"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('E48').Value = '  +3.01%  '
$ws.Range('E50').Value = '  +4.79%  '
$ws.Range('E51').Value = '  +8.76%  '
